$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.057.92"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "1.667.09"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'215.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "'0.5102"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.2669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.06393"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'21.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "'0.07440"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "1.672.25"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "'4.516"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'0.5809"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "'0.000008511"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "'64.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "26.005.65"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "'4.919"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D20").Value = "'10.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "'189.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").Value = "'6.191"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'144.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "'7.613"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E26").Value = "  +4.13%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'0.06625"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.77%  "
$ws.Range("D29").Value = "'1.330"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").Value = "'1.311"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").Value = "'3.549"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'3.511"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "'1.017"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").Value = "'0.6148"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").Value = "'2.369"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "'2.687"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").Value = "'6.391"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.50%  "
$ws.Range("D39").Value = "1.093.79"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").Value = "'0.01593"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").Value = "'0.8687"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").Value = "'101.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").Value = "1.813.96"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "'56.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "'8.095"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'0.4287"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  +3.24%  "
